$wb = $excel.ActiveWorkbook

# Sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").ClearContents()
$ws.Range("N17").Value = 0
$ws.Range("H43").Value = 2421
$ws.Range("I43").Value = 2494.5
$ws.Range("J43").Value = 2347.5
$ws.Range("K43").Value = 2494.5
$ws.Range("L43").Value = 2347.5
$ws.Range("M43").Value = -2425.5
$ws.Range("N43").Value = -2485.5
$ws.Range("H51").Value = 3450
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 3450
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 3450
$ws.Range("N51").Value = -4418
$ws.Range("H86").Value = 5760.773
$ws.Range("I86").Value = 4998.7856
$ws.Range("J86").Value = 7094.25
$ws.Range("K86").Value = 4998.7856
$ws.Range("L86").Value = 7094.25
$ws.Range("M86").Value = -3875.7856
$ws.Range("N86").Value = -9340.25
$ws.Range("H89").Value = 5760.773
$ws.Range("I89").Value = 4998.7856
$ws.Range("J89").Value = 7094.25
$ws.Range("K89").Value = 24993.928
$ws.Range("L89").Value = 35471.25
$ws.Range("M89").Value = -19377.928
$ws.Range("N89").Value = -46703.25
$ws.Range("H101").Value = 11111845
$ws.Range("I101").Value = 33333580
$ws.Range("J101").Value = 977.3333
$ws.Range("K101").Value = 100000740
$ws.Range("L101").Value = 2931.9999
$ws.Range("M101").Value = -99999118
$ws.Range("N101").Value = -6175.9999
$ws.Range("H118").Value = 1836
$ws.Range("I118").Value = 1836
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 5508
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -3851
$ws.Range("H137").Value = 1461
$ws.Range("I137").Value = 674
$ws.Range("J137").Value = 2248
$ws.Range("K137").Value = 2022
$ws.Range("L137").Value = 6744
$ws.Range("M137").Value = 528
$ws.Range("N137").Value = -11844
$ws.Range("H138").Value = 3700.7349
$ws.Range("I138").Value = 3587.5881
$ws.Range("J138").Value = 3729.879
$ws.Range("K138").Value = 10762.7643
$ws.Range("L138").Value = 11189.637
$ws.Range("M138").Value = -5622.764299999999
$ws.Range("N138").Value = -21469.637
$ws.Range("H141").Value = 4777.8335
$ws.Range("I141").Value = 4343.5
$ws.Range("J141").Value = 6949.5
$ws.Range("K141").Value = 13030.5
$ws.Range("L141").Value = 20848.5
$ws.Range("M141").Value = -7850.5
$ws.Range("N141").Value = -31208.5

# Sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 120
$ws.Range("I15").Value = 120
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 120
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 230
$ws.Range("H32").Value = 5183.269
$ws.Range("I32").Value = 4615.25
$ws.Range("J32").Value = 11999.5
$ws.Range("K32").Value = 4615.25
$ws.Range("L32").Value = 11999.5
$ws.Range("M32").Value = -4328.25
$ws.Range("N32").Value = -12573.5
$ws.Range("H45").Value = 1900
$ws.Range("I45").Value = 1900
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1900
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1523
$ws.Range("H61").Value = 3318.3215
$ws.Range("I61").Value = 1468.7222
$ws.Range("J61").Value = 6647.6
$ws.Range("K61").Value = 1468.7222
$ws.Range("L61").Value = 6647.6
$ws.Range("M61").Value = -1256.7222
$ws.Range("N61").Value = -7071.6
$ws.Range("H74").Value = 2391.6
$ws.Range("I74").Value = 2013.7059
$ws.Range("J74").Value = 4533
$ws.Range("K74").Value = 2013.7059
$ws.Range("L74").Value = 4533
$ws.Range("M74").Value = -1139.7059
$ws.Range("N74").Value = -6281
$ws.Range("H77").Value = 2391.6
$ws.Range("I77").Value = 2013.7059
$ws.Range("J77").Value = 4533
$ws.Range("K77").Value = 10068.5295
$ws.Range("L77").Value = 22665
$ws.Range("M77").Value = -5700.529500000001
$ws.Range("N77").Value = -31401
$ws.Range("H127").Value = 149999
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 149999
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 149999
$ws.Range("N127").Value = -159919
$ws.Range("H132").Value = 2904
$ws.Range("I132").Value = 2904
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8712
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -6182
$ws.Range("H136").Value = 3318.3215
$ws.Range("I136").Value = 1468.7222
$ws.Range("J136").Value = 6647.6
$ws.Range("K136").Value = 4406.1666
$ws.Range("L136").Value = 19942.8
$ws.Range("M136").Value = -1856.1666
$ws.Range("N136").Value = -25042.8

# Sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1818.8
$ws.Range("I20").Value = 1818.8
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 1818.8
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -1571.8
$ws.Range("H134").Value = 800
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 800
$ws.Range("K134").Value = 0
$ws.Range("L134").ClearContents()
$ws.Range("M134").Value = 2400
$ws.Range("N134").Value = -7470

# Sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3016
$ws.Range("I16").Value = 3850
$ws.Range("J16").Value = 1904
$ws.Range("K16").Value = 3850
$ws.Range("L16").Value = 1904
$ws.Range("M16").Value = -3563
$ws.Range("N16").Value = -2478
$ws.Range("H58").Value = 4838
$ws.Range("I58").Value = 4500
$ws.Range("J58").Value = 5007
$ws.Range("K58").Value = 4500
$ws.Range("L58").Value = 5007
$ws.Range("M58").Value = -4297
$ws.Range("N58").Value = -5413
$ws.Range("H113").Value = 3016
$ws.Range("I113").Value = 3850
$ws.Range("J113").Value = 1904
$ws.Range("K113").Value = 3850
$ws.Range("L113").Value = 1904
$ws.Range("M113").Value = -1680
$ws.Range("N113").Value = -6244
$ws.Range("H129").Value = 115775.2
$ws.Range("I129").Value = 19000
$ws.Range("J129").Value = 139969
$ws.Range("K129").Value = 19000
$ws.Range("L129").Value = 139969
$ws.Range("M129").Value = -14000
$ws.Range("N129").Value = -149969
$ws.Range("H136").Value = 4838
$ws.Range("I136").Value = 4500
$ws.Range("J136").Value = 5007
$ws.Range("K136").Value = 13500
$ws.Range("L136").Value = 15021
$ws.Range("M136").Value = -10950
$ws.Range("N136").Value = -20121

# Sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 170158.77
$ws.Range("I2").Value = 366678.66
$ws.Range("J2").Value = 111202.8
$ws.Range("K2").Value = 2200071.96
$ws.Range("L2").Value = 667216.8
$ws.Range("M2").Value = -2199958.96
$ws.Range("N2").Value = -667442.8
$ws.Range("H121").Value = 12707.3125
$ws.Range("I121").Value = 30827
$ws.Range("J121").Value = 6667.4165
$ws.Range("K121").Value = 92481
$ws.Range("L121").Value = 20002.2495
$ws.Range("M121").Value = -91171
$ws.Range("N121").Value = -22622.2495
$ws.Range("H122").Value = 667.5
$ws.Range("I122").Value = 688.6667
$ws.Range("J122").Value = 646.3333
$ws.Range("K122").Value = 6198.0003
$ws.Range("L122").Value = 5816.9997
$ws.Range("M122").Value = -3748.0003
$ws.Range("N122").Value = -10716.9997

# Sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 1124.75
$ws.Range("I22").Value = 999.5
$ws.Range("J22").Value = 1250
$ws.Range("K22").Value = 999.5
$ws.Range("L22").Value = 1250
$ws.Range("M22").Value = -470.5
$ws.Range("N22").Value = -2308
$ws.Range("H70").Value = 5012
$ws.Range("I70").Value = 4000
$ws.Range("J70").Value = 6024
$ws.Range("K70").Value = 4000
$ws.Range("L70").Value = 6024
$ws.Range("M70").Value = -3730
$ws.Range("N70").Value = -6564
$ws.Range("H73").Value = 5012
$ws.Range("I73").Value = 4000
$ws.Range("J73").Value = 6024
$ws.Range("K73").Value = 4000
$ws.Range("L73").Value = 6024
$ws.Range("M73").Value = -3064
$ws.Range("N73").Value = -7896
$ws.Range("H80").Value = 9732.4375
$ws.Range("I80").Value = 2998
$ws.Range("J80").Value = 10694.5
$ws.Range("K80").Value = 2998
$ws.Range("L80").Value = 10694.5
$ws.Range("M80").Value = -2000
$ws.Range("N80").Value = -12690.5
$ws.Range("H83").Value = 9732.4375
$ws.Range("I83").Value = 2998
$ws.Range("J83").Value = 10694.5
$ws.Range("K83").Value = 14990
$ws.Range("L83").Value = 53472.5
$ws.Range("M83").Value = -9998
$ws.Range("N83").Value = -63456.5
$ws.Range("H107").Value = 307.92856
$ws.Range("I107").Value = 254.6923
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 254.6923
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 1665.3077
$ws.Range("N107").Value = -4840
$ws.Range("H113").Value = 2199.3333
$ws.Range("I113").Value = 2199.3333
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2199.3333
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -29.33329999999978
$ws.Range("H129").Value = 116519.336
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 116519.336
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 116519.336
$ws.Range("N129").Value = -126519.336
$ws.Range("H132").Value = 2170.8572
$ws.Range("I132").Value = 2039.4
$ws.Range("J132").Value = 2499.5
$ws.Range("K132").Value = 6118.200000000001
$ws.Range("L132").Value = 7498.5
$ws.Range("M132").Value = -3588.200000000001
$ws.Range("N132").Value = -12558.5

# Sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3651.8096
$ws.Range("I40").Value = 3234.6667
$ws.Range("J40").Value = 4208
$ws.Range("K40").Value = 3234.6667
$ws.Range("L40").Value = 4208
$ws.Range("M40").Value = -3098.6667
$ws.Range("N40").Value = -4480
$ws.Range("H82").Value = 1094.5
$ws.Range("I82").Value = 1126.3334
$ws.Range("J82").Value = 999
$ws.Range("K82").Value = 1126.3334
$ws.Range("L82").Value = 999
$ws.Range("M82").Value = -765.3334
$ws.Range("N82").Value = -1721
$ws.Range("H85").Value = 1094.5
$ws.Range("I85").Value = 1126.3334
$ws.Range("J85").Value = 999
$ws.Range("K85").Value = 1126.3334
$ws.Range("L85").Value = 999
$ws.Range("M85").Value = 121.6666
$ws.Range("N85").Value = -3495
$ws.Range("H93").Value = 937.25
$ws.Range("I93").Value = 975
$ws.Range("J93").Value = 899.5
$ws.Range("K93").Value = 975
$ws.Range("L93").Value = 899.5
$ws.Range("M93").Value = 273
$ws.Range("N93").Value = -3395.5
$ws.Range("H128").Value = 54999
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 54999
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 54999
$ws.Range("N128").Value = -64959
$ws.Range("H132").Value = 3534.5
$ws.Range("I132").Value = 2719.125
$ws.Range("J132").Value = 4621.6665
$ws.Range("K132").Value = 8157.375
$ws.Range("L132").Value = 13864.9995
$ws.Range("M132").Value = -5627.375
$ws.Range("N132").Value = -18924.9995
$ws.Range("H140").Value = 74999.5
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 74999.5
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 74999.5
$ws.Range("N140").Value = -85359.5

# Sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 10000
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 10000
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 10000
$ws.Range("N22").Value = -10586
$ws.Range("H130").Value = 64685
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 64685
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 64685
$ws.Range("N130").Value = -74725
$ws.Range("H136").Value = 2663.3333
$ws.Range("I136").Value = 2663.3333
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 7989.999899999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -5439.999899999999
